# CommonNames_required_edited.xlsx — refresh the taxonomic-assignment table.
#
# The sheet went from 5 assigned/unassigned species rows down to just the
# "Unassigned" row plus a single new record: Microstomus kitt (Lemon sole).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Mareca americana" (bird, row 5) and "Myrophis vafer" (eel, row 6)
# rows first so row numbers for the earlier rows don't shift while we work.
$ws.Range("A5:J6").EntireRow.Delete()

# Drop the old "Fundulus heteroclitus or majalis" (row 2) and "Cololabis saira"
# (row 3) rows. This leaves row 2 = "Unassigned" (previously row 4) and row 3 empty.
$ws.Range("A2:J3").EntireRow.Delete()

# Populate row 3 with the new species' full taxonomy.
$newSpecies = @(
    "Microstomus kitt",
    "Lemon sole",
    "Teleost Fish",
    "Animalia",
    "Chordata",
    "Teleostei",
    "Carangiformes",
    "Pleuronectidae",
    "Microstomus",
    "kitt"
)
for ($col = 1; $col -le $newSpecies.Length; $col++) {
    $ws.Cells.Item(3, $col).Value = $newSpecies[$col - 1]
}

# Re-fit the columns to the now-shorter content.
$ws.Columns.Item(1).ColumnWidth = 15.5546875 - 5/6
$ws.Columns.Item(2).ColumnWidth = 14.6640625 - 5/6
$ws.Columns.Item(3).ColumnWidth = 14 - 5/6
for ($col = 4; $col -le 10; $col++) {
    $ws.Columns.Item($col).ColumnWidth = 10.77734375 - 5/6
}

# Leave the selection where it ended up after the edit (first empty row).
$ws.Range("J4").Select()
